$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 20:22"

# Update country names (column A) for rows whose rank/order changed
$ws.Range("A61").Value = "Marruecos"
$ws.Range("A62").Value = "Croacia"
$ws.Range("A70").Value = "Barein"
$ws.Range("A71").Value = "Lituania"
$ws.Range("A112").Value = "Guinea"
$ws.Range("A113").Value = "Georgia"
$ws.Range("A114").Value = "Isla de Man"
$ws.Range("A115").Value = "Consejo Danes para los Refugiados"
$ws.Range("A195").Value = "Republica de Africa Central"
$ws.Range("A196").Value = "Sierra Leona"
$ws.Range("A197").Value = "Nicaragua"
$ws.Range("A198").Value = "Islas Turcas y Caicos"
$ws.Range("A206").Value = "Sahara Occidental"
$ws.Range("A207").Value = "Santo Tome y Principe"

# Update numeric statistics (columns B-H) for rows with updated data
$ws.Range("B4").Value = 521816
$ws.Range("C4").Value = 18940
$ws.Range("D4").Value = 28587
$ws.Range("E4").Value = 473160
$ws.Range("F4").Value = 10961
$ws.Range("G4").Value = 1322
$ws.Range("H4").Value = 20069

$ws.Range("B17").Value = 20173
$ws.Range("C17").Value = 384
$ws.Range("E17").Value = 18916
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 1084

$ws.Range("B24").Value = 8439
$ws.Range("C24").Value = 839
$ws.Range("E24").Value = 7182

$ws.Range("B61").Value = 1545
$ws.Range("C61").Value = 97
$ws.Range("D61").Value = 146
$ws.Range("E61").Value = 1288
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 111

$ws.Range("B62").Value = 1534
$ws.Range("C62").Value = 39
$ws.Range("D62").Value = 323
$ws.Range("E62").Value = 1190
$ws.Range("F62").Value = 32
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 21

$ws.Range("B70").Value = 1040
$ws.Range("C70").Value = 115
$ws.Range("D70").Value = 555
$ws.Range("E70").Value = 479
$ws.Range("F70").Value = 3
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 6

$ws.Range("B71").Value = 1026
$ws.Range("C71").Value = 27
$ws.Range("D71").Value = 54
$ws.Range("E71").Value = 949
$ws.Range("F71").Value = 14
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 23

$ws.Range("B112").Value = 250
$ws.Range("C112").Value = 38
$ws.Range("D112").Value = 17
$ws.Range("E112").Value = 233
$ws.Range("F112").Value = 0
$ws.Range("H112").Value = 0

$ws.Range("B113").Value = 242
$ws.Range("C113").Value = 8
$ws.Range("D113").Value = 60
$ws.Range("E113").Value = 179
$ws.Range("F113").Value = 6
$ws.Range("H113").Value = 3

$ws.Range("B114").Value = 226
$ws.Range("C114").Value = 25
$ws.Range("D114").Value = 112
$ws.Range("E114").Value = 113
$ws.Range("F114").Value = 11
$ws.Range("H114").Value = 1

$ws.Range("B115").Value = 223
$ws.Range("D115").Value = 16
$ws.Range("E115").Value = 187
$ws.Range("H115").Value = 20

$ws.Range("F116").Value = 1

$ws.Range("C197").Value = 1

$ws.Range("C198").Value = 0
